$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 80, shifting existing rows 80-126 down to 81-127
$ws.Rows.Item(80).Insert()

# Fill in the new row 80 with its data (mirrors the format of surrounding rows)
$ws.Cells.Item(80, 1).Value = 11
$ws.Cells.Item(80, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(80, 3).Value = "Bíobío"
$ws.Cells.Item(80, 4).Value = 44806
$ws.Cells.Item(80, 5).Value = 8
$ws.Cells.Item(80, 6).Value = "Fruta"
$ws.Cells.Item(80, 7).Value = 100108
$ws.Cells.Item(80, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(80, 9).Value = 100108002
$ws.Cells.Item(80, 10).Value = "Mango"
$ws.Cells.Item(80, 11).Value = "Sin especificar"
$ws.Cells.Item(80, 12).Value = "Primera"
$ws.Cells.Item(80, 13).Value = 200
$ws.Cells.Item(80, 14).Value = 9000
$ws.Cells.Item(80, 15).Value = 9500
$ws.Cells.Item(80, 16).Value = 9250
$ws.Cells.Item(80, 17).Value = '$/bandeja 4 kilos'
$ws.Cells.Item(80, 18).Value = "Brasil"
$ws.Cells.Item(80, 19).Value = 2312
$ws.Cells.Item(80, 20).Value = 4

# Match the date cell number format used by the other date cells in column D
$ws.Cells.Item(80, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
